# Updated cryptos list - apply Price (D) and Volume(1h) (E) changes per diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.708.19"
$ws.Range("E2").Value = "  +1.41%  "

$ws.Range("D3").Value = "'1.644.90"
$ws.Range("E3").Value = "  -0.39%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").Value = "'213.47"
$ws.Range("E5").Value = "  +0.18%  "

$ws.Range("D6").Value = "'0.534"
$ws.Range("E6").Value = "  +3.89%  "

$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("D8").Value = "'23.23"
$ws.Range("E8").Value = "  -1.04%  "

$ws.Range("D9").Value = "'0.260"
$ws.Range("E9").Value = "  +0.34%  "

$ws.Range("E10").Value = "  +0.09%  "

$ws.Range("D11").Value = "'0.0892"
$ws.Range("E11").Value = "  +1.67%  "

$ws.Range("D12").Value = "'1.877.80"
$ws.Range("E12").Value = "  -0.41%  "

$ws.Range("D13").Value = "'1.653.70"
$ws.Range("E13").Value = "  -0.09%  "

$ws.Range("E14").Value = "  -0.83%  "

$ws.Range("D15").Value = "'0.564"
$ws.Range("E15").Value = "  -0.94%  "

$ws.Range("D16").Value = "'64.28"
$ws.Range("E16").Value = "  -1.77%  "

$ws.Range("D17").Value = "'27.682.76"
$ws.Range("E17").Value = "  +1.27%  "

$ws.Range("D18").Value = "'231.72"
$ws.Range("E18").Value = "  +0.34%  "

$ws.Range("E19").Value = "  +0.01%  "

$ws.Range("D20").Value = "'7.68"
$ws.Range("E20").Value = "  +3.73%  "

$ws.Range("E21").Value = "  -0.01%  "

$ws.Range("D22").Value = "'4.32"
$ws.Range("E22").Value = "  -1.02%  "

$ws.Range("D23").Value = "'10.10"
$ws.Range("E23").Value = "  +8.67%  "

$ws.Range("E24").Value = "  -3.79%  "

$ws.Range("D25").Value = "'149.99"
$ws.Range("E25").Value = "  +1.80%  "

$ws.Range("E26").Value = "  -1.80%  "

$ws.Range("D27").Value = "'0.112"
$ws.Range("E27").Value = "  +0.56%  "

$ws.Range("E28").Value = "  -0.08%  "

$ws.Range("E29").Value = "  -0.94%  "

$ws.Range("E30").Value = "  -0.43%  "

$ws.Range("E31").Value = "  -1.94%  "

$ws.Range("E32").Value = "  +0.51%  "

$ws.Range("D33").Value = "'1.444.73"
$ws.Range("E33").Value = "  +0.95%  "

$ws.Range("E34").Value = "  +0.76%  "

$ws.Range("E35").Value = "  +1.65%  "

$ws.Range("D36").Value = "'2.35"
$ws.Range("E36").Value = "  -1.00%  "

$ws.Range("E37").Value = "  -0.21%  "

$ws.Range("D38").Value = "'0.885"
$ws.Range("E38").Value = "  -2.21%  "

$ws.Range("E39").Value = "  -0.54%  "

$ws.Range("D40").Value = "'0.899"
$ws.Range("E40").Value = "  +14.10%  "

$ws.Range("E41").Value = "  -1.79%  "

$ws.Range("E42").Value = "  +0.01%  "

$ws.Range("D43").Value = "'5.71"
$ws.Range("E43").Value = "  +2.35%  "

$ws.Range("E44").Value = "  -0.59%  "

$ws.Range("D45").Value = "'66.21"
$ws.Range("E45").Value = "  +1.92%  "

$ws.Range("E46").Value = "  +1.76%  "

$ws.Range("D47").Value = "'1.787.07"
$ws.Range("E47").Value = "  -0.39%  "

$ws.Range("E48").Value = "  +2.48%  "

$ws.Range("D49").Value = "'86.55"
$ws.Range("E49").Value = "  -1.71%  "

$ws.Range("E50").Value = "  +1.58%  "

$ws.Range("E51").Value = "  -2.10%  "

